$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now describes the remaining student (MONTBULEAU--GENTELET Titouan)
$ws.Range("A2").Value = "MONTBULEAU--GENTELET"
$ws.Range("B2").Value = "TITOUAN"
$ws.Range("G2").Value = "stage"
$ws.Range("H2").Value = "MAAT PHARMA"
$ws.Range("I2").Value = "Lyon"
$ws.Range("J2").Value = "BERGÉ"

# Phone numbers must stay text (leading zero) rather than become a number
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "0617421317"
$ws.Range("K2").Style = "Normal"

$ws.Range("L2").Value = "URRUTY"

# The other students (rows 3-6) are no longer part of the list
$ws.Rows("3:6").Delete()
